{"js": "// Finish API6: re-highlight the checkUserName / checkRegister /\n// searchFreeMeetingRoom API entries from yellow to green, and move the\n// (Word-managed) \"_GoBack\" bookmark so it sits after the last edited\n// entry, i.e. right after the \"searchFreeMeetingRoom\" run.\n\nconst body = context.document.body;\n\nasync function highlightGreen(apiName) {\n  // Match the leading \"/\" together with the API name so both runs\n  // (\"/\" and the identifier) get re-highlighted, same as the source edit.\n  const results = body.search(\"/\" + apiName, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const range = results.items[0];\n  // \"BrightGreen\" round-trips to OOXML <w:highlight w:val=\"green\"/>.\n  range.font.highlightColor = \"BrightGreen\";\n  await context.sync();\n  return range;\n}\n\nawait highlightGreen(\"checkUserName\");\nawait highlightGreen(\"checkRegister\");\nconst lastRange = await highlightGreen(\"searchFreeMeetingRoom\");\n\n// \"_GoBack\" is Word's own last-edit-location bookmark: drop the old one\n// and drop a fresh one, collapsed, right after the last edited run.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst endRange = lastRange.getRange(\"End\");\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Finish API6: re-highlight the checkUserName / checkRegister /\n# searchFreeMeetingRoom API entries from yellow to green, and move the\n# (Word-managed) \"_GoBack\" bookmark so it again sits after the last edit\n# location, i.e. right after the \"searchFreeMeetingRoom\" run.\n\n$d = $word.ActiveDocument\n\nfunction Set-ApiHighlightGreen($apiName) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $apiName\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    if ($find.Execute()) {\n        # wdBrightGreen (4) round-trips to OOXML <w:highlight w:val=\"green\"/>\n        $rng.HighlightColorIndex = 4\n    }\n    return $rng\n}\n\nSet-ApiHighlightGreen(\"checkUserName\") | Out-Null\nSet-ApiHighlightGreen(\"checkRegister\") | Out-Null\n$lastRng = Set-ApiHighlightGreen(\"searchFreeMeetingRoom\")\n\n# Re-adding the bookmark under this name moves Word's special \"_GoBack\"\n# bookmark here (and removes it from wherever it previously was).\n$d.Bookmarks.Add(\"_GoBack\", $lastRng) | Out-Null\n"}
